$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("I8").Value = 1.87

# Row 10
$ws.Range("G10").Value = 1.87

# Row 11
$ws.Range("G11").Value = 3.25
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 2.18
$ws.Range("M11").Value = 2.4
$ws.Range("N11").Value = 2.27
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 1.47
$ws.Range("Q11").Value = 2.32
$ws.Range("R11").Value = 2.02
$ws.Range("S11").Value = 1.62
$ws.Range("T11").Value = 7.7
$ws.Range("U11").Value = 15.5
$ws.Range("V11").Value = 12.5
$ws.Range("W11").Value = 45
$ws.Range("X11").Value = 35
$ws.Range("Y11").Value = 55
$ws.Range("Z11").Value = 7
$ws.Range("AA11").Value = 6.1
$ws.Range("AB11").Value = 18.5
$ws.Range("AC11").Value = 120
$ws.Range("AE11").Value = 6
$ws.Range("AF11").Value = 9.25
$ws.Range("AG11").Value = 9.5
$ws.Range("AH11").Value = 20
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 40

# Row 13
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9
$ws.Range("N13").Value = 2.2
$ws.Range("O13").Value = 1.62

# Row 14
$ws.Range("G14").Value = 7.9
$ws.Range("H14").Value = 3.95
$ws.Range("I14").Value = 1.37
$ws.Range("L14").Value = 1.27
$ws.Range("M14").Value = 3.4
$ws.Range("N14").Value = 1.82
$ws.Range("O14").Value = 1.8
$ws.Range("R14").Value = 2.12
$ws.Range("S14").Value = 1.64
$ws.Range("T14").Value = 14.5
$ws.Range("V14").Value = 20
$ws.Range("X14").Value = 80
$ws.Range("Y14").Value = 70
$ws.Range("Z14").Value = 9.25
$ws.Range("AB14").Value = 17.5
$ws.Range("AC14").Value = 90
$ws.Range("AD14").Value = 500
$ws.Range("AE14").Value = 5.1
$ws.Range("AF14").Value = 5.1
$ws.Range("AG14").Value = 7
$ws.Range("AH14").Value = 7.3
$ws.Range("AI14").Value = 10
$ws.Range("AJ14").Value = 24

# Row 15
$ws.Range("G15").Value = 4.65
$ws.Range("H15").Value = 3.55
$ws.Range("L15").Value = 1.27
$ws.Range("M15").Value = 3.4
$ws.Range("N15").Value = 1.8
$ws.Range("O15").Value = 1.82
$ws.Range("R15").Value = 1.87
$ws.Range("S15").Value = 1.84
$ws.Range("X15").Value = 35
$ws.Range("Z15").Value = 10
$ws.Range("AA15").Value = 6.1
$ws.Range("AE15").Value = 5.8
$ws.Range("AF15").Value = 6.5
$ws.Range("AG15").Value = 6.9
$ws.Range("AI15").Value = 10.75

# Row 16
$ws.Range("G16").Value = 1.78
$ws.Range("H16").Value = 3.5
$ws.Range("I16").Value = 3.85
$ws.Range("L16").Value = 1.31
$ws.Range("M16").Value = 3.15
$ws.Range("N16").Value = 1.88
$ws.Range("P16").Value = 1.42
$ws.Range("Q16").Value = 2.65
$ws.Range("R16").Value = 1.9
$ws.Range("S16").Value = 1.81
$ws.Range("T16").Value = 5.8
$ws.Range("U16").Value = 6.9
$ws.Range("V16").Value = 7.1
$ws.Range("W16").Value = 11.5
$ws.Range("X16").Value = 12
$ws.Range("Y16").Value = 22
$ws.Range("Z16").Value = 9.5
$ws.Range("AA16").Value = 6
$ws.Range("AB16").Value = 13.5
$ws.Range("AC16").Value = 60
$ws.Range("AE16").Value = 9
$ws.Range("AF16").Value = 16.5
$ws.Range("AG16").Value = 11
$ws.Range("AH16").Value = 45
$ws.Range("AI16").Value = 29
$ws.Range("AJ16").Value = 35

# Row 18
$ws.Range("G18").Value = 5.6
$ws.Range("H18").Value = 4.75
$ws.Range("I18").Value = 1.45
$ws.Range("J18").Value = 1.03
$ws.Range("K18").Value = 9.75
$ws.Range("L18").Value = 1.14
$ws.Range("M18").Value = 4.9
$ws.Range("N18").Value = 1.44
$ws.Range("O18").Value = 2.6
$ws.Range("Q18").Value = 3.55
$ws.Range("T18").Value = 22
$ws.Range("U18").Value = 40
$ws.Range("V18").Value = 18
$ws.Range("W18").Value = 100
$ws.Range("X18").Value = 45
$ws.Range("Z18").Value = 9.75
$ws.Range("AA18").Value = 9.75
$ws.Range("AB18").Value = 15.5
$ws.Range("AE18").Value = 10.25
$ws.Range("AH18").Value = 11
$ws.Range("AI18").Value = 10.5

# Row 19
$ws.Range("G19").Value = 1.95
$ws.Range("I19").Value = 3.4
$ws.Range("M19").Value = 3.55
$ws.Range("S19").Value = 2.12
$ws.Range("T19").Value = 9
$ws.Range("U19").Value = 10.75
$ws.Range("V19").Value = 8.5
$ws.Range("W19").Value = 18
$ws.Range("Y19").Value = 21
$ws.Range("AE19").Value = 12
$ws.Range("AF19").Value = 20
$ws.Range("AG19").Value = 11.75
$ws.Range("AH19").Value = 45
$ws.Range("AI19").Value = 28

# Row 20
$ws.Range("G20").Value = 1.17
$ws.Range("I20").Value = 11.75
$ws.Range("N20").Value = 1.47
$ws.Range("O20").Value = 2.5
$ws.Range("R20").Value = 2.05
$ws.Range("S20").Value = 1.69
$ws.Range("T20").Value = 7.6
$ws.Range("U20").Value = 5.7
$ws.Range("V20").Value = 8.5
$ws.Range("W20").Value = 5.9
$ws.Range("X20").Value = 8.75
$ws.Range("Y20").Value = 24
$ws.Range("Z20").Value = 17
$ws.Range("AA20").Value = 11.25
$ws.Range("AB20").Value = 23
$ws.Range("AC20").Value = 90
$ws.Range("AD20").Value = 450
$ws.Range("AE20").Value = 28
$ws.Range("AF20").Value = 80
$ws.Range("AG20").Value = 32
$ws.Range("AH20").Value = 300
$ws.Range("AI20").Value = 120
$ws.Range("AJ20").Value = 90

# Row 21
$ws.Range("G21").Value = 2.4
$ws.Range("H21").Value = 3.2
$ws.Range("I21").Value = 2.75
$ws.Range("L21").Value = 1.31
$ws.Range("M21").Value = 2.87
$ws.Range("R21").Value = 1.72
$ws.Range("S21").Value = 1.88
$ws.Range("T21").Value = 8
$ws.Range("U21").Value = 11.75
$ws.Range("W21").Value = 25
$ws.Range("X21").Value = 20
$ws.Range("AA21").Value = 6.2
$ws.Range("AB21").Value = 14
$ws.Range("AC21").Value = 65
$ws.Range("AD21").Value = 500
$ws.Range("AE21").Value = 8.25
$ws.Range("AF21").Value = 13.5
$ws.Range("AG21").Value = 10.25
$ws.Range("AI21").Value = 24

# Row 24
$ws.Range("G24").Value = 3.25
$ws.Range("H24").Value = 3.7
$ws.Range("I24").Value = 2
$ws.Range("U24").Value = 21

# Row 25
$ws.Range("G25").Value = 3.1
$ws.Range("H25").Value = 3.6
$ws.Range("J25").Value = 1.06
$ws.Range("K25").Value = 8
$ws.Range("L25").Value = 1.36
$ws.Range("M25").Value = 3
$ws.Range("N25").Value = 2.1
$ws.Range("O25").Value = 1.7
$ws.Range("P25").Value = 1.41
$ws.Range("Q25").Value = 2.62
$ws.Range("R25").Value = 1.91
$ws.Range("S25").Value = 1.8
$ws.Range("X25").Value = 26
$ws.Range("Z25").Value = 9
$ws.Range("AB25").Value = 17
$ws.Range("AC25").Value = 51
$ws.Range("AD25").Value = 700
$ws.Range("AE25").Value = 7
$ws.Range("AF25").Value = 10
$ws.Range("AI25").Value = 19

# Row 26
$ws.Range("P26").Value = 1.37

# Row 29
$ws.Range("N29").Value = 1.13

# Row 30
$ws.Range("N30").Value = 1.72
$ws.Range("O30").Value = 2.05

# Row 32
$ws.Range("G32").Value = 2.45
$ws.Range("I32").Value = 2.75
$ws.Range("K32").Value = 12
$ws.Range("N32").Value = 1.77
$ws.Range("T32").Value = 9.5
$ws.Range("V32").Value = 9.5
$ws.Range("W32").Value = 23
$ws.Range("AH32").Value = 29

# Row 33
$ws.Range("O33").Value = 1.47
